$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new yellow highlight style (fills count 3->4, cellXfs 3->4) ---
$ws.Range("E12").Interior.Color = 65535

# --- Column width changes ---
$ws.Columns.Item(3).ColumnWidth = 69.16666666666667
$ws.Columns.Item(4).ColumnWidth = 39.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 59.166666666666664

# --- Row data updates (rows 2-13) ---
# Row 2
$ws.Cells.Item(2,1).Value = "'1326670"
$ws.Cells.Item(2,2).Value = "https://aiesec.org/opportunity/global-talent/1326670"
$ws.Cells.Item(2,3).Value = "TIM Operations Assistant Intern"
$ws.Cells.Item(2,4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(2,5).Value = "No"
$ws.Cells.Item(2,6).Value = "1 applicant"
$ws.Cells.Item(2,7).Value = "6 - 18 Months"
$ws.Cells.Item(2,8).Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"

# Row 3
$ws.Cells.Item(3,1).Value = "'1326669"
$ws.Cells.Item(3,2).Value = "https://aiesec.org/opportunity/global-talent/1326669"
$ws.Cells.Item(3,3).Value = "EB Sales & Supply Chain Management Assistant"
$ws.Cells.Item(3,4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(3,5).Value = "No"
$ws.Cells.Item(3,6).Value = "1 applicant"
$ws.Cells.Item(3,7).Value = "6 - 18 Months"
$ws.Cells.Item(3,8).Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"

# Row 4
$ws.Cells.Item(4,1).Value = "'1326664"
$ws.Cells.Item(4,2).Value = "https://aiesec.org/opportunity/global-talent/1326664"
$ws.Cells.Item(4,3).Value = "Market Research & Insights Coordinator Trainee ONLY EU"
$ws.Cells.Item(4,4).Value = "Bruxelles, Belgio"
$ws.Cells.Item(4,5).Value = "No"
$ws.Cells.Item(4,6).Value = "5 applicants"
$ws.Cells.Item(4,7).Value = "6 - 18 Months"
$ws.Cells.Item(4,8).Value = "UCB"

# Row 5
$ws.Cells.Item(5,1).Value = "'1326636"
$ws.Cells.Item(5,2).Value = "https://aiesec.org/opportunity/global-talent/1326636"
$ws.Cells.Item(5,3).Value = "Computer and AI Coordinator"
$ws.Cells.Item(5,4).Value = "London, UK"
$ws.Cells.Item(5,5).Value = "No"
$ws.Cells.Item(5,6).Value = "49 applicants"
$ws.Cells.Item(5,7).Value = "3 - 6 Months"
$ws.Cells.Item(5,8).Value = "Capital Care Homes"

# Row 6
$ws.Cells.Item(6,1).Value = "'1326291"
$ws.Cells.Item(6,2).Value = "https://aiesec.org/opportunity/global-talent/1326291"
$ws.Cells.Item(6,3).Value = "Marketing Trainee Health Systems"
$ws.Cells.Item(6,4).Value = "Santiago, Región Metropolitana, Chile"
$ws.Cells.Item(6,5).Value = "No"
$ws.Cells.Item(6,6).Value = "3 applicants"
$ws.Cells.Item(6,7).Value = "6 - 18 Months"
$ws.Cells.Item(6,8).Value = "Philips Chilena"

# Row 7
$ws.Cells.Item(7,1).Value = "'1326162"
$ws.Cells.Item(7,2).Value = "https://aiesec.org/opportunity/global-talent/1326162"
$ws.Cells.Item(7,3).Value = "Tech Sales Development Representative( swedish Only)"
$ws.Cells.Item(7,4).Value = "Bournemouth, Royaume-Uni"
$ws.Cells.Item(7,5).Value = "No"
$ws.Cells.Item(7,6).Value = "7 applicants"
$ws.Cells.Item(7,7).Value = "6 - 18 Months"
$ws.Cells.Item(7,8).Value = "EIMS Ltd"

# Row 8
$ws.Cells.Item(8,1).Value = "'1326160"
$ws.Cells.Item(8,2).Value = "https://aiesec.org/opportunity/global-talent/1326160"
$ws.Cells.Item(8,3).Value = "Tech Sales Development Representative( dutch  Only)"
$ws.Cells.Item(8,4).Value = "Bournemouth, Royaume-Uni"
$ws.Cells.Item(8,5).Value = "No"
$ws.Cells.Item(8,6).Value = "3 applicants"
$ws.Cells.Item(8,7).Value = "6 - 18 Months"
$ws.Cells.Item(8,8).Value = "EIMS Ltd"

# Row 9
$ws.Cells.Item(9,1).Value = "'1326159"
$ws.Cells.Item(9,2).Value = "https://aiesec.org/opportunity/global-talent/1326159"
$ws.Cells.Item(9,3).Value = "Tech Sales Development Representative( Spanish Only)"
$ws.Cells.Item(9,4).Value = "Bournemouth, Royaume-Uni"
$ws.Cells.Item(9,5).Value = "No"
$ws.Cells.Item(9,6).Value = "33 applicants"
$ws.Cells.Item(9,7).Value = "6 - 18 Months"
$ws.Cells.Item(9,8).Value = "EIMS Ltd"

# Row 10
$ws.Cells.Item(10,1).Value = "'1326156"
$ws.Cells.Item(10,2).Value = "https://aiesec.org/opportunity/global-talent/1326156"
$ws.Cells.Item(10,3).Value = "Tech Sales Development Representative( French/ Swiss/ Belgian Only)"
$ws.Cells.Item(10,4).Value = "Bournemouth, Royaume-Uni"
$ws.Cells.Item(10,5).Value = "No"
$ws.Cells.Item(10,6).Value = "19 applicants"
$ws.Cells.Item(10,7).Value = "6 - 18 Months"
$ws.Cells.Item(10,8).Value = "EIMS Ltd"

# Row 11
$ws.Cells.Item(11,1).Value = "'1326152"
$ws.Cells.Item(11,2).Value = "https://aiesec.org/opportunity/global-talent/1326152"
$ws.Cells.Item(11,3).Value = "Tech Sales Development Representative( German / Austrian Only)"
$ws.Cells.Item(11,4).Value = "Bournemouth, Royaume-Uni"
$ws.Cells.Item(11,5).Value = "No"
$ws.Cells.Item(11,6).Value = "3 applicants"
$ws.Cells.Item(11,7).Value = "6 - 18 Months"
$ws.Cells.Item(11,8).Value = "EIMS Ltd"

# Row 12
$ws.Cells.Item(12,1).Value = "'1324011"
$ws.Cells.Item(12,2).Value = "https://aiesec.org/opportunity/global-talent/1324011"
$ws.Cells.Item(12,3).Value = "ACE Program | Quality Engineer (Thai)"
$ws.Cells.Item(12,4).Value = "Hyderabad, Telangana, India"
$ws.Cells.Item(12,5).Value = "Yes"
$ws.Cells.Item(12,6).Value = "13 applicants"
$ws.Cells.Item(12,7).Value = "6 - 18 Months"
$ws.Cells.Item(12,8).Value = "Tata Consultancy Services Ltd."

# Row 13
$ws.Cells.Item(13,1).Value = "'1320966"
$ws.Cells.Item(13,2).Value = "https://aiesec.org/opportunity/global-talent/1320966"
$ws.Cells.Item(13,3).Value = "Sales and Marketing Intern"
$ws.Cells.Item(13,4).Value = "Mumbai, Maharashtra, India"
$ws.Cells.Item(13,5).Value = "No"
$ws.Cells.Item(13,6).Value = "15 applicants"
$ws.Cells.Item(13,7).Value = "6 - 18 Months"
$ws.Cells.Item(13,8).Value = "Agrocel Industries Private Limited"
